$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.629732666666667
$ws.Range("H2").Value = 7.889198
$ws.Range("I2").Value = 0.07156737804735891
$ws.Range("J2").Value = 0.07156737804735891
$ws.Range("M2").Value = 70.23436
$ws.Range("N2").Value = 210.70308
$ws.Range("O2").Value = 0.7023186840741513
$ws.Range("P2").Value = 0.7023186840741513
$ws.Range("Q2").Value = 184.6975908144267
$ws.Range("R2").Value = 1662.27831732984
$ws.Range("S2").Value = 0.05026310677285842
$ws.Range("T2").Value = 0.05026310677285842
$ws.Range("G3").Value = 2.629732666666667
$ws.Range("H3").Value = 7.889198
$ws.Range("I3").Value = 0.07156737804735891
$ws.Range("J3").Value = 0.07156737804735891
$ws.Range("O3").Value = 0.06551129587759326
$ws.Range("P3").Value = 0.06551129587759325
$ws.Range("Q3").Value = 17.22833066255867
$ws.Range("R3").Value = 155.054975963028
$ws.Range("S3").Value = 0.004688471678444102
$ws.Range("T3").Value = 0.004688471678444101
$ws.Range("G4").Value = 2.629732666666667
$ws.Range("H4").Value = 7.889198
$ws.Range("I4").Value = 0.07156737804735891
$ws.Range("J4").Value = 0.07156737804735891
$ws.Range("M4").Value = 4.268944666666666
$ws.Range("N4").Value = 12.806834
$ws.Range("O4").Value = 0.04268793224112385
$ws.Range("P4").Value = 0.04268793224112385
$ws.Range("Q4").Value = 11.22618324212578
$ws.Range("R4").Value = 101.035649179132
$ws.Range("S4").Value = 0.003055063384760552
$ws.Range("T4").Value = 0.003055063384760552
$ws.Range("G5").Value = 2.629732666666667
$ws.Range("H5").Value = 7.889198
$ws.Range("I5").Value = 0.07156737804735891
$ws.Range("J5").Value = 0.07156737804735891
$ws.Range("M5").Value = 18.948881
$ws.Range("N5").Value = 56.846643
$ws.Range("O5").Value = 0.1894820878071316
$ws.Range("P5").Value = 0.1894820878071315
$ws.Range("Q5").Value = 49.83049136247934
$ws.Range("R5").Value = 448.474422262314
$ws.Range("S5").Value = 0.01356073621129584
$ws.Range("T5").Value = 0.01356073621129584
$ws.Range("I6").Value = 0.493312042610523
$ws.Range("J6").Value = 0.493312042610523
$ws.Range("M6").Value = 70.23436
$ws.Range("N6").Value = 210.70308
$ws.Range("O6").Value = 0.7023186840741513
$ws.Range("P6").Value = 0.7023186840741513
$ws.Range("Q6").Value = 1273.115604844627
$ws.Range("R6").Value = 11458.04044360164
$ws.Range("S6").Value = 0.3464622646041542
$ws.Range("T6").Value = 0.3464622646041542
$ws.Range("I7").Value = 0.493312042610523
$ws.Range("J7").Value = 0.493312042610523
$ws.Range("O7").Value = 0.06551129587759326
$ws.Range("P7").Value = 0.06551129587759325
$ws.Range("S7").Value = 0.03231751118343787
$ws.Range("T7").Value = 0.03231751118343786
$ws.Range("I8").Value = 0.493312042610523
$ws.Range("J8").Value = 0.493312042610523
$ws.Range("M8").Value = 4.268944666666666
$ws.Range("N8").Value = 12.806834
$ws.Range("O8").Value = 0.04268793224112385
$ws.Range("P8").Value = 0.04268793224112385
$ws.Range("Q8").Value = 77.38178394950245
$ws.Range("R8").Value = 696.436055545522
$ws.Range("S8").Value = 0.02105847104868841
$ws.Range("T8").Value = 0.02105847104868841
$ws.Range("I9").Value = 0.493312042610523
$ws.Range("J9").Value = 0.493312042610523
$ws.Range("M9").Value = 18.948881
$ws.Range("N9").Value = 56.846643
$ws.Range("O9").Value = 0.1894820878071316
$ws.Range("P9").Value = 0.1894820878071315
$ws.Range("Q9").Value = 343.4802580310244
$ws.Range("R9").Value = 3091.322322279219
$ws.Range("S9").Value = 0.09347379577424256
$ws.Range("T9").Value = 0.09347379577424254
$ws.Range("G10").Value = 7.550656333333333
$ws.Range("H10").Value = 22.651969
$ws.Range("I10").Value = 0.2054888252189962
$ws.Range("J10").Value = 0.2054888252189962
$ws.Range("M10").Value = 70.23436
$ws.Range("N10").Value = 210.70308
$ws.Range("O10").Value = 0.7023186840741513
$ws.Range("P10").Value = 0.7023186840741513
$ws.Range("Q10").Value = 530.3155151516133
$ws.Range("R10").Value = 4772.839636364521
$ws.Range("S10").Value = 0.1443186413197487
$ws.Range("T10").Value = 0.1443186413197487
$ws.Range("G11").Value = 7.550656333333333
$ws.Range("H11").Value = 22.651969
$ws.Range("I11").Value = 0.2054888252189962
$ws.Range("J11").Value = 0.2054888252189962
$ws.Range("O11").Value = 0.06551129587759326
$ws.Range("P11").Value = 0.06551129587759325
$ws.Range("Q11").Value = 49.46708297725934
$ws.Range("R11").Value = 445.203746795334
$ws.Range("S11").Value = 0.01346183922846071
$ws.Range("T11").Value = 0.01346183922846071
$ws.Range("G12").Value = 7.550656333333333
$ws.Range("H12").Value = 22.651969
$ws.Range("I12").Value = 0.2054888252189962
$ws.Range("J12").Value = 0.2054888252189962
$ws.Range("M12").Value = 4.268944666666666
$ws.Range("N12").Value = 12.806834
$ws.Range("O12").Value = 0.04268793224112385
$ws.Range("P12").Value = 0.04268793224112385
$ws.Range("Q12").Value = 32.23333408401622
$ws.Range("R12").Value = 290.100006756146
$ws.Range("S12").Value = 0.008771893047256653
$ws.Range("T12").Value = 0.008771893047256653
$ws.Range("G13").Value = 7.550656333333333
$ws.Range("H13").Value = 22.651969
$ws.Range("I13").Value = 0.2054888252189962
$ws.Range("J13").Value = 0.2054888252189962
$ws.Range("M13").Value = 18.948881
$ws.Range("N13").Value = 56.846643
$ws.Range("O13").Value = 0.1894820878071316
$ws.Range("P13").Value = 0.1894820878071315
$ws.Range("Q13").Value = 143.0764883322297
$ws.Range("R13").Value = 1287.688394990067
$ws.Range("S13").Value = 0.03893645162353015
$ws.Range("T13").Value = 0.03893645162353015
$ws.Range("G14").Value = 8.437784666666667
$ws.Range("H14").Value = 25.313354
$ws.Range("I14").Value = 0.2296317541231219
$ws.Range("J14").Value = 0.2296317541231219
$ws.Range("M14").Value = 70.23436
$ws.Range("N14").Value = 210.70308
$ws.Range("O14").Value = 0.7023186840741513
$ws.Range("P14").Value = 0.7023186840741513
$ws.Range("Q14").Value = 592.6224058811467
$ws.Range("R14").Value = 5333.60165293032
$ws.Range("S14").Value = 0.16127467137739
$ws.Range("T14").Value = 0.16127467137739
$ws.Range("G15").Value = 8.437784666666667
$ws.Range("H15").Value = 25.313354
$ws.Range("I15").Value = 0.2296317541231219
$ws.Range("J15").Value = 0.2296317541231219
$ws.Range("O15").Value = 0.06551129587759326
$ws.Range("P15").Value = 0.06551129587759325
$ws.Range("Q15").Value = 55.27898182938267
$ws.Range("R15").Value = 497.510836464444
$ws.Range("S15").Value = 0.01504347378725058
$ws.Range("T15").Value = 0.01504347378725058
$ws.Range("G16").Value = 8.437784666666667
$ws.Range("H16").Value = 25.313354
$ws.Range("I16").Value = 0.2296317541231219
$ws.Range("J16").Value = 0.2296317541231219
$ws.Range("M16").Value = 4.268944666666666
$ws.Range("N16").Value = 12.806834
$ws.Range("O16").Value = 0.04268793224112385
$ws.Range("P16").Value = 0.04268793224112385
$ws.Range("Q16").Value = 36.02043585124844
$ws.Range("R16").Value = 324.183922661236
$ws.Range("S16").Value = 0.009802504760418239
$ws.Range("T16").Value = 0.009802504760418239
$ws.Range("G17").Value = 8.437784666666667
$ws.Range("H17").Value = 25.313354
$ws.Range("I17").Value = 0.2296317541231219
$ws.Range("J17").Value = 0.2296317541231219
$ws.Range("M17").Value = 18.948881
$ws.Range("N17").Value = 56.846643
$ws.Range("O17").Value = 0.1894820878071316
$ws.Range("P17").Value = 0.1894820878071315
$ws.Range("Q17").Value = 159.8865775522914
$ws.Range("R17").Value = 1438.979197970622
$ws.Range("S17").Value = 0.04351110419806303
$ws.Range("T17").Value = 0.04351110419806303
